$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 4.5
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("W2").Value = 5.5
$ws.Range("AD2").Value = 5.5
$ws.Range("AF2").Value = 81
$ws.Range("AO2").Value = 15
$ws.Range("G5").Value = 1.55
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 2.1
$ws.Range("L5").Value = 6
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.93
$ws.Range("AE5").Value = 19
$ws.Range("AI5").Value = 19
$ws.Range("AL5").Value = 51
$ws.Range("AY5").Value = 41
$ws.Range("AZ5").Value = 126
$ws.Range("BA5").Value = 151
$ws.Range("BB5").Value = 301
$ws.Range("Q6").Value = 1.9
$ws.Range("R6").Value = 1.95
$ws.Range("H7").Value = 2.87
$ws.Range("J7").Value = 2.8
$ws.Range("K7").Value = 1.95
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.52
$ws.Range("W7").Value = 6.2
$ws.Range("AA7").Value = 21
$ws.Range("AC7").Value = 7.1
$ws.Range("AE7").Value = 14.5
$ws.Range("AL7").Value = 40
$ws.Range("AQ7").Value = 50
$ws.Range("AR7").Value = 90
$ws.Range("G9").Value = 2.38
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 2.9
$ws.Range("J9").Value = 3.2
$ws.Range("L9").Value = 3.75
$ws.Range("W9").Value = 7
$ws.Range("X9").Value = 11
$ws.Range("Z9").Value = 23
$ws.Range("AA9").Value = 21
$ws.Range("AC9").Value = 7.5
$ws.Range("AE9").Value = 17
$ws.Range("AH9").Value = 15
$ws.Range("AI9").Value = 12
$ws.Range("AJ9").Value = 34
$ws.Range("AK9").Value = 29
$ws.Range("AL9").Value = 41
$ws.Range("AN9").Value = 4.33
$ws.Range("AV9").Value = 67
$ws.Range("BB9").Value = 251
$ws.Range("G12").Value = 1.29
$ws.Range("I12").Value = 11
$ws.Range("L12").Value = 9.5
$ws.Range("U12").Value = 2.5
$ws.Range("V12").Value = 1.5
$ws.Range("Y12").Value = 9.5
$ws.Range("Z12").Value = 7.5
$ws.Range("AC12").Value = 10
$ws.Range("AG12").Value = 19
$ws.Range("AI12").Value = 29
$ws.Range("AK12").Value = 81
$ws.Range("AL12").Value = 81
$ws.Range("AN12").Value = 3.1
$ws.Range("AR12").Value = 41
$ws.Range("AU12").Value = 11
$ws.Range("AW12").Value = 10
$ws.Range("AZ12").Value = 301
$ws.Range("BA12").Value = 301
